# HW2: add a third column (C) of computed values next to the existing
# student-name (A) / score (B) columns, and move the active selection to C1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cValues = @(
    88,84,57,43,87,92,76.6666666666667,77.095238095238102,77.523809523809504,
    77.952380952381006,78.380952380952394,78.809523809523796,79.238095238095298,
    79.6666666666667,80.095238095238102,80.523809523809504,80.952380952381006,
    81.380952380952394,81.809523809523796,82.238095238095298,82.6666666666667,
    83.095238095238102,83.523809523809504,83.952380952381006,84.380952380952394,
    84.809523809523796,85.238095238095298,85.6666666666667,86.095238095238102,
    86.523809523809504,86.952380952381006,87.380952380952394,87.809523809523796,
    88.238095238095298,88.6666666666667,89.095238095238102,89.523809523809604,
    89.952380952381006,90.380952380952394,90.809523809523796,91.238095238095298,
    91.6666666666667,92.095238095238102,92.523809523809604,92.952380952381006,
    93.380952380952394,93.809523809523796,94.238095238095298,94.6666666666667,
    95.095238095238102
)

for ($i = 0; $i -lt $cValues.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 3).Value = $cValues[$i]
}

$ws.Range("C1").Select()
